# Generate Report for Handoff
# The "ad1f4133-009a-4d78-8668-3180561983d5" entry has finished translation
# and is now ready for handoff; update its status and handoff timestamps
# across the Overview sheet and each per-language sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: row for ad1f4133-...-md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("D3").Value = "2016-36-19 06:36:44"

# --- zh-cn sheet: row for ad1f4133-...-md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("E3").Value = "2016-03-19 06:36:42"

# --- de-de sheet: row for ad1f4133-...-md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("E3").Value = "2016-03-19 06:36:44"
